# Update cryptos list with latest pricing/volume figures (scheduled refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> (new Price text or $null if unchanged, new Volume(1h) text)
$updates = @(
    @{ Row = 2;  D = "42.924.62"; E = "  +0.34%  " },
    @{ Row = 3;  D = "2.290.74";  E = "  +1.75%  " },
    @{ Row = 4;  D = $null;       E = "  -0.12%  " },
    @{ Row = 5;  D = "252.12";    E = "  +0.88%  " },
    @{ Row = 6;  D = "0.649";     E = "  +4.13%  " },
    @{ Row = 7;  D = "75.57";     E = "  +6.99%  " },
    @{ Row = 9;  D = $null;       E = "  -3.60%  " },
    @{ Row = 10; D = "39.33";     E = "  -0.50%  " },
    @{ Row = 11; D = "0.0977";    E = "  +0.82%  " },
    @{ Row = 12; D = $null;       E = "  +0.61%  " },
    @{ Row = 13; D = "0.107";     E = "  +1.65%  " },
    @{ Row = 14; D = "2.635.23";  E = "  +1.82%  " },
    @{ Row = 15; D = "15.05";     E = "  +2.01%  " },
    @{ Row = 16; D = $null;       E = "  -1.41%  " },
    @{ Row = 17; D = "2.293.81";  E = "  +2.27%  " },
    @{ Row = 18; D = "42.811.15"; E = "  +0.26%  " },
    @{ Row = 19; D = $null;       E = "  +1.46%  " },
    @{ Row = 20; D = $null;       E = "  -0.52%  " },
    @{ Row = 21; D = "72.35";     E = "  -0.83%  " },
    @{ Row = 22; D = "237.43";    E = "  +1.67%  " },
    @{ Row = 23; D = "2.16";      E = "  +4.56%  " },
    @{ Row = 24; D = $null;       E = "  -1.12%  " },
    @{ Row = 25; D = $null;       E = "  -0.12%  " },
    @{ Row = 26; D = "11.31";     E = "  -1.96%  " },
    @{ Row = 27; D = $null;       E = "  -0.91%  " },
    @{ Row = 28; D = "2.12";      E = "  -3.60%  " },
    @{ Row = 29; D = "167.47";    E = "  +0.01%  " },
    @{ Row = 30; D = "21.05";     E = "  +0.49%  " },
    @{ Row = 31; D = $null;       E = "  +9.37%  " },
    @{ Row = 32; D = $null;       E = "  -4.58%  " },
    @{ Row = 33; D = "0.126";     E = "  -1.29%  " },
    @{ Row = 34; D = "31.06";     E = "  +0.38%  " },
    @{ Row = 35; D = $null;       E = "  +1.85%  " },
    @{ Row = 36; D = "4.60";      E = "  +4.76%  " },
    @{ Row = 37; D = $null;       E = "  +1.81%  " },
    @{ Row = 38; D = $null;       E = "  -5.15%  " },
    @{ Row = 39; D = "13.60";     E = "  +8.10%  " },
    @{ Row = 40; D = $null;       E = "  -0.61%  " },
    @{ Row = 41; D = "5.94";      E = "  +1.81%  " },
    @{ Row = 42; D = $null;       E = "  +4.25%  " },
    @{ Row = 43; D = $null;       E = "  +1.64%  " },
    @{ Row = 44; D = "61.17";     E = "  -2.43%  " },
    @{ Row = 45; D = "4.85";      E = "  +0.05%  " },
    @{ Row = 46; D = "105.32";    E = "  +11.08%  " },
    @{ Row = 47; D = "0.101";     E = "  -1.62%  " },
    @{ Row = 48; D = $null;       E = "  -0.25%  " },
    @{ Row = 49; D = "1.16";      E = "  -0.33%  " },
    @{ Row = 50; D = $null;       E = "  -1.70%  " },
    @{ Row = 51; D = "4.23";      E = "  -1.76%  " }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($null -ne $u.D) {
        $cell = $ws.Cells.Item($r, 4)
        # These "Price" column entries are plain text labels in the sheet
        # (e.g. "252.12", "0.649"). Force text so Excel doesn't silently
        # reinterpret them as numbers (which would also eat trailing zeros,
        # e.g. "4.60" -> 4.6).
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
    }
    $ws.Cells.Item($r, 5).Value = $u.E
}
